$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26: amend existing Indie Project entry, bump hours, grow row height ---
$ws.Range("B26").Value = 5.5
$ws.Range("D26").Value = "Indie Project: setting up test database, properties, copying Database class, creating sql for cleaning database before testing, starting UserDaoTest, setting up log4j properties, changing mysql username and pw, trying to delete sensitive files from GitHub, troubleshooting why my cleandb.sql isn't actually running..."
$ws.Rows.Item(26).RowHeight = 45

# --- Clear the three old entries that are being removed / relocated ---
$ws.Range("D28").Clear()
$ws.Range("D32").Clear()
$ws.Range("D37").Clear()

# --- New row 27: next day's entry ---
$ws.Range("A27").Value = 43520
$ws.Range("A27").NumberFormat = "d-mmm"
$ws.Range("D27").Value = "Week 4ish: Getting passwords out of repo for Week 1 exercise with intention to figure out why this project reads cleandb.sql and it doesn't get read in my project tests"
$ws.Rows.Item(27).RowHeight = 30

# --- New row 31 ---
$ws.Range("D31").Value = "Sun AM - 1 hr"

# --- Row 33 now holds a different note (old content relocates to row 36) ---
$ws.Range("D33").Value = "NEXT: remove logs from github -- they currently only have the removed user/password"

# --- Row 35 now holds the relocated "NOT RESOLVED" note ---
$ws.Range("D35").Value = "NOT RESOLVED: getting correct path for mysqldump (don't need to do it yet so defering…)"

# --- New row 36: relocated "Now I have a path..." note ---
$ws.Range("D36").Value = "Now I have a path to get at mysqldump but I have an access problem for writing the dump to the locations I choose"

# --- New row 38: relocated log4j note ---
$ws.Range("D38").Value = "I have a log4j problem about renaming with the date etc, appears to be the same in week 1 exercise and project."

# --- New row 40 ---
$ws.Range("D40").Value = "NB: new password so files/setup need to change in prior repos"

# --- New row 42 ---
$ws.Range("D42").Value = "WHY won't cleandb.sql run in my project?"

# --- Update the active selection to match the latest edit location ---
$ws.Range("D48").Select()
